$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K-SVMeans")

# Fill in the new column C values (k = 3 data) for rows 13-17
$ws.Range("C13").Value = 58.23
$ws.Range("C14").Value = 55.67
$ws.Range("C15").Value = 50.26
$ws.Range("C16").Value = 54.89
$ws.Range("C17").Value = 53.24

# Add the AVERAGE formula in C18, matching the style of B18
$ws.Range("C18").Formula = "=AVERAGE(C13:C17)"
$ws.Range("C18").Font.Color = $ws.Range("B18").Font.Color

# Update the selected cell to D17 as in the diff
$ws.Range("D17").Select()
